$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (column D) and Volume(1h) (column E) cells to match latest scrape.
# Values are written with a leading apostrophe to force text interpretation
# (several look like numbers, e.g. "1.00", "0.539"), then the cell style is
# reset to "Normal" so no extra quote-prefix / text-format style is introduced.

$ws.Range("D2").Value = '''63.093.39'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '''  -0.70%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''3.044.01'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '''  -0.90%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = '''1.00'
$ws.Range("D4").Style = "Normal"

$ws.Range("D5").Value = '''588.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '''  -0.31%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''151.68'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '''  -0.94%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = '''  -0.02%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''0.539'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '''  -1.67%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''3.045.68'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '''  -1.01%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = '''  -1.85%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = '''5.80'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '''  -1.03%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '''0.449'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '''  -2.92%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = '''  -2.90%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''36.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '''  -2.52%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("E15").Value = '''  +1.41%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''3.542.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '''  -1.07%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''7.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '''  -1.30%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''63.085.28'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '''  -0.67%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = '''3.040.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '''  -1.16%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = '''477.73'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '''  +0.33%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''14.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '''  -3.11%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''0.706'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '''  -1.95%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("E23").Value = '''  -0.65%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("E24").Value = '''  +1.36%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''81.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '''  +0.53%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = '''12.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '''  -2.91%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = '''10.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '''  +7.10%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = '''1.00'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '''  +0.25%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = '''7.37'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '''  +0.20%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("E30").Value = '''  -0.58%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = '''1.00'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '''  -0.01%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("E32").Value = '''  -0.18%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = '''27.66'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '''  +1.05%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = '''  -3.11%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("E35").Value = '''  +0.64%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = '''  -4.67%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = '''  -3.22%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = '''  -3.77%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("E39").Value = '''  +0.08%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''9.27'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '''  -0.69%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''50.46'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '''  -0.04%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = '''435.44'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '''  -2.97%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = '''  +0.13%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("E44").Value = '''  +1.98%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = '''0.0362'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '''  -0.73%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = '''2.824.35'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '''  +0.42%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = '''38.34'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '''  -5.13%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = '''128.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '''  -1.91%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D50").Value = '''25.11'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '''  -0.50%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = '''2.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '''  -2.31%  '
$ws.Range("E51").Style = "Normal"
